# Updated cryptos list (GitHub Actions style price/volume refresh).
# Cells are plain text (prices use "." as both thousands & decimal
# separators e.g. "66.633.36", so they must never be coerced to a real
# number by Excel's auto-detection). We force text by stamping the
# cell's NumberFormat to "@" before writing the value, then clear the
# format back off afterwards so no stray style index is left behind
# (matches original cells, which carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '66.633.36'
Set-TextValue 'E2' '  +0.16%  '
Set-TextValue 'D3' '3.527.87'
Set-TextValue 'E3' '  -1.91%  '
Set-TextValue 'D5' '607.43'
Set-TextValue 'E5' '  -0.11%  '
Set-TextValue 'D6' '143.35'
Set-TextValue 'D7' '3.526.43'
Set-TextValue 'E7' '  -1.93%  '
Set-TextValue 'E8' '  -0.13%  '
Set-TextValue 'E9' '  +4.21%  '
Set-TextValue 'D10' '7.72'
Set-TextValue 'E10' '  -3.72%  '
Set-TextValue 'E11' '  -4.44%  '
Set-TextValue 'D12' '0.408'
Set-TextValue 'E12' '  -1.84%  '
Set-TextValue 'D13' '4.120.83'
Set-TextValue 'E13' '  -2.02%  '
Set-TextValue 'E14' '  -6.29%  '
Set-TextValue 'D15' '28.67'
Set-TextValue 'E15' '  -3.84%  '
Set-TextValue 'D16' '3.520.24'
Set-TextValue 'E16' '  -1.79%  '
Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '66.513.18'
Set-TextValue 'E17' '  -0.10%  '
Set-TextValue 'B18' 'TRON'
Set-TextValue 'C18' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D18' '0.117'
Set-TextValue 'E18' '  +0.39%  '
Set-TextValue 'D19' '10.78'
Set-TextValue 'E19' '  -6.87%  '
Set-TextValue 'E20' '  -3.62%  '
Set-TextValue 'D21' '14.61'
Set-TextValue 'E21' '  -3.30%  '
Set-TextValue 'D22' '423.22'
Set-TextValue 'E22' '  -0.99%  '
Set-TextValue 'E23' '  -5.03%  '
Set-TextValue 'D24' '77.08'
Set-TextValue 'E24' '  -2.10%  '
Set-TextValue 'D25' '3.673.74'
Set-TextValue 'E25' '  -1.81%  '
Set-TextValue 'E26' '  +0.00%  '
Set-TextValue 'E27' '  -5.48%  '
Set-TextValue 'D28' '7.91'
Set-TextValue 'E28' '  -4.91%  '
Set-TextValue 'E29' '  -1.92%  '
Set-TextValue 'D30' '8.94'
Set-TextValue 'E30' '  -5.07%  '
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  +0.12%  '
Set-TextValue 'D32' '3.533.14'
Set-TextValue 'E32' '  -1.69%  '
Set-TextValue 'E33' '  -2.36%  '
Set-TextValue 'D34' '24.17'
Set-TextValue 'E34' '  -4.90%  '
Set-TextValue 'E35' '  +0.03%  '
Set-TextValue 'D36' '1.32'
Set-TextValue 'E36' '  -9.64%  '
Set-TextValue 'E37' '  -3.59%  '
Set-TextValue 'E38' '  -4.01%  '
Set-TextValue 'D39' '173.66'
Set-TextValue 'E39' '  -2.26%  '
Set-TextValue 'E40' '  -7.64%  '
Set-TextValue 'D41' '0.0812'
Set-TextValue 'E41' '  -5.11%  '
Set-TextValue 'D42' '4.99'
Set-TextValue 'E42' '  -4.79%  '
Set-TextValue 'D43' '0.852'
Set-TextValue 'E43' '  -5.04%  '
Set-TextValue 'D44' '45.50'
Set-TextValue 'D45' '1.78'
Set-TextValue 'E45' '  -6.42%  '
Set-TextValue 'D46' '0.999'
Set-TextValue 'E46' '  +0.01%  '
Set-TextValue 'E47' '  -7.83%  '
Set-TextValue 'D48' '7.07'
Set-TextValue 'E48' '  -1.85%  '
Set-TextValue 'E49' '  -4.39%  '
Set-TextValue 'D50' '22.80'
Set-TextValue 'E50' '  -4.79%  '
Set-TextValue 'D51' '0.905'
Set-TextValue 'E51' '  -5.20%  '
